# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (with the quarterly fund-holding detail)
# right after the "总计" (totals) summary sheet, and updates the totals
# sheet with the new quarter's summary row (shifting the older quarters
# down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet, positioned right after "总计"
#    and therefore right before the (old) first quarter sheet "2022-Q3".
# ---------------------------------------------------------------------
$totalsSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$q4Sheet = $wb.Worksheets.Add($null, $totalsSheet)
$q4Sheet.Name = "2022-Q4"

# Header row (bold, centered, thin-bordered - matches the other quarter sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Fund detail rows. D/E/F/G columns are textual percentages / figures
# (stored as text, same as on the other quarter sheets), H is numeric.
$data = @(
    @("010583", "富国蓝筹精选股票（QDII）美元",       "13.38", "91.97", "6.14", "0.8215", 2),
    @("007455", "富国蓝筹精选股票（QDII）人民币",     "13.38", "91.97", "6.14", "0.8215", 2),
    @("000934", "国富大中华精选混合（QDII）",          "20.61", "87.91", "3.02", "0.6224", 9),
    @("006370", "国富大中华精选混合（QDII）美元",      "20.61", "87.91", "3.02", "0.6224", 9),
    @("010671", "景顺长城大中华混合（QDII）美元A",     "12.44", "86.89", "4.21", "0.5237", 8),
    @("262001", "景顺长城大中华混合（QDII）人民币A",   "12.42", "86.89", "4.21", "0.5229", 8),
    @("100055", "富国全球科技互联网股票（QDII）",      "3.86",  "94.32", "5.53", "0.2135", 3),
    @("016988", "景顺长城大中华混合（QDII）人民币C",   "0.02",  "86.89", "4.21", "0.0008", 8)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $q4Sheet.Cells.Item($row, 1).Value = $i
    $q4Sheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $q4Sheet.Cells.Item($row, 3).Value = "'" + $rec[1]
    $q4Sheet.Cells.Item($row, 4).Value = "'" + $rec[2]
    $q4Sheet.Cells.Item($row, 5).Value = "'" + $rec[3]
    $q4Sheet.Cells.Item($row, 6).Value = "'" + $rec[4]
    $q4Sheet.Cells.Item($row, 7).Value = "'" + $rec[5]
    $q4Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# Match the header / index-column styling used throughout the workbook:
# bold font, centered horizontally, top-aligned vertically, thin border.
$styledRanges = @($q4Sheet.Range("B1:H1"), $q4Sheet.Range("A2:A9"))
foreach ($rng in $styledRanges) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

$q4Sheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add the 2022-Q4 row at the top of
#    the data (row 2) and push the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = @(
    @("2022-Q4", 8, 4.15),
    @("2022-Q3", 8, 3.98),
    @("2022-Q2", 3, 1.52),
    @("2022-Q1", 4, 1.77),
    @("2021-Q4", 3, 1.8)
)

for ($i = 0; $i -lt $summary.Count; $i++) {
    $row = $i + 2
    $rec = $summary[$i]

    $totalsSheet.Cells.Item($row, 1).Value = $i
    $totalsSheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $totalsSheet.Cells.Item($row, 3).Value = $rec[1]
    $totalsSheet.Cells.Item($row, 4).Value = $rec[2]
}
